$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.975.54"
$ws.Range("E2").Value = "  -1.42%  "

$ws.Range("D3").Value = "2.452.57"
$ws.Range("E3").Value = "  -3.69%  "

$ws.Range("E4").Value = "  +0.15%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "524.41"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "129.94"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.32%  "

$ws.Range("E7").Value = "  +0.26%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.565"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.12%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0976"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.41%  "

$ws.Range("E10").Value = "  -2.27%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "4.96"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.62%  "

$ws.Range("E12").Value = "  -4.04%  "

$ws.Range("D13").Value = "2.890.60"
$ws.Range("E13").Value = "  -3.53%  "

$ws.Range("D14").Value = "57.932.23"
$ws.Range("E14").Value = "  -1.52%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.58"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.88%  "

$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("D17").Value = "2.457.14"
$ws.Range("E17").Value = "  -3.53%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "10.38"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.61%  "

$ws.Range("E19").Value = "  -2.04%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "311.55"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.04%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.13"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("E24").Value = "  -2.54%  "

$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("D26").Value = "2.566.50"
$ws.Range("E26").Value = "  -3.53%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.156"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.92%  "

$ws.Range("E28").Value = "  -3.08%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "174.79"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.61%  "

$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").Value = "  -3.17%  "

$ws.Range("E31").Value = "  -2.78%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.57%  "

$ws.Range("E33").Value = "  -6.17%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.84"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.80%  "

$ws.Range("E37").Value = "  -7.16%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.78"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -5.04%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "36.35"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.23%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.806"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("E41").Value = "  -4.55%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.39"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.13%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.584"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.38%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "258.22"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -8.25%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.99%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "124.31"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -7.45%  "

$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("E48").Value = "  -3.09%  "

$ws.Range("E49").Value = "  -3.11%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "17.06"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.78%  "

$ws.Range("E51").Value = "  -5.20%  "
